$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DE relative area: DE (row 9) as a percentage of MB (row 8)
$ws.Range("F9").Formula = "=C9/C8*100"
$ws.Range("F9").NumberFormat = "0.0000"

# LF relative area: LF (row 10) as a ratio of FM_Total (row 2)
$ws.Range("F10").Formula = "=C10/C2"
$ws.Range("F10").Style = "Normal"

$ws.Range("F10").Select() | Out-Null
